# Generate Report for Handoff
# - Refreshes the "Latest Handoff Datetime" / "Latest HO Xliff Generate Date"
#   timestamps for the files that were just (re)handed off.
# - Marks those same rows' Priority column as "ht" (handoff type) on the
#   per-language detail sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$rows = @(7, 9, 11, 12, 13, 14)

# Refresh the handoff timestamps that are shared across the just-regenerated rows.
foreach ($r in $rows) {
    $overview.Range("G$r").Value = "2016-09-01 00:23:11"
    $zhcn.Range("H$r").Value     = "2016-09-01 00:23:02"
    $dede.Range("H$r").Value     = "2016-09-01 00:23:11"
}

# Flag these rows as handoff-type priority on both language sheets.
foreach ($r in $rows) {
    $zhcn.Range("E$r").Value = "ht"
    $dede.Range("E$r").Value = "ht"
}
